$wb = $excel.ActiveWorkbook

# --- choices sheet: insert new choice row for "re" right after the last
# existing ce5ds76 (National Red List Category) choice, i.e. before row 315,
# shifting every subsequent choices row down by one.
$choices = $wb.Worksheets.Item("choices")
$choices.Activate()
$choices.Rows.Item(315).Insert()

$choices.Cells.Item(315, 1).Value = "ce5ds76"
$choices.Cells.Item(315, 2).Value = "re"
$choices.Cells.Item(315, 3).Value = "Regionally extinct (RE). Select only for National Redlists."

$choices.Range("D304").Select()

# --- survey sheet: column A gained an explicit width and the view scrolled to
# + selected the "regional_redlist" (select_one ce5ds76 / National Red List
# Category) question row while making this edit.
$survey = $wb.Worksheets.Item("survey")
$survey.Columns.Item(1).ColumnWidth = 16.3333333333
$survey.Activate()
$survey.Range("A663:C663").Select()

Write-Host "done"
